$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.676.06"
$ws.Range("E2").Value = "  -1.06%  "

# Row 3
$ws.Range("D3").Value = "3.790.72"
$ws.Range("E3").Value = "  -0.46%  "

# Row 5
$ws.Range("D5").Value = "'596.21"
$ws.Range("E5").Value = "  +0.10%  "

# Row 6
$ws.Range("D6").Value = "'166.76"
$ws.Range("E6").Value = "  -0.99%  "

# Row 7
$ws.Range("D7").Value = "3.789.77"
$ws.Range("E7").Value = "  -0.53%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 10
$ws.Range("E10").Value = "  -1.03%  "

# Row 11
$ws.Range("D11").Value = "'6.35"
$ws.Range("E11").Value = "  -1.13%  "

# Row 12
$ws.Range("E12").Value = "  -0.52%  "

# Row 13
$ws.Range("E13").Value = "  -2.98%  "

# Row 14
$ws.Range("D14").Value = "'36.07"
$ws.Range("E14").Value = "  -0.53%  "

# Row 15
$ws.Range("D15").Value = "4.426.19"
$ws.Range("E15").Value = "  -0.45%  "

# Row 16
$ws.Range("D16").Value = "3.831.68"
$ws.Range("E16").Value = "  +0.69%  "

# Row 17
$ws.Range("D17").Value = "'18.59"
$ws.Range("E17").Value = "  +3.06%  "

# Row 18
$ws.Range("D18").Value = "67.677.60"
$ws.Range("E18").Value = "  -1.09%  "

# Row 19
$ws.Range("E19").Value = "  +1.01%  "

# Row 20
$ws.Range("E20").Value = "  +0.10%  "

# Row 21
$ws.Range("E21").Value = "  -9.58%  "

# Row 22
$ws.Range("D22").Value = "'459.85"
$ws.Range("E22").Value = "  -1.54%  "

# Row 23
$ws.Range("D23").Value = "'0.699"
$ws.Range("E23").Value = "  -0.39%  "

# Row 24
$ws.Range("D24").Value = "'0.0000153"
$ws.Range("E24").Value = "  +1.55%  "

# Row 25
$ws.Range("D25").Value = "'83.36"
$ws.Range("E25").Value = "  -1.05%  "

# Row 26
$ws.Range("D26").Value = "'12.04"
$ws.Range("E26").Value = "  +0.77%  "

# Row 27
$ws.Range("E27").Value = "  -3.75%  "

# Row 28
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.08%  "

# Row 29
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'10.01"
$ws.Range("E29").Value = "  -1.52%  "

# Row 30
$ws.Range("D30").Value = "3.939.75"
$ws.Range("E30").Value = "  -0.43%  "

# Row 31
$ws.Range("E31").Value = "  -0.05%  "

# Row 32
$ws.Range("D32").Value = "'2.25"
$ws.Range("E32").Value = "  +3.34%  "

# Row 33
$ws.Range("D33").Value = "'7.23"
$ws.Range("E33").Value = "  -1.47%  "

# Row 34
$ws.Range("D34").Value = "'29.61"
$ws.Range("E34").Value = "  -2.11%  "

# Row 35
$ws.Range("E35").Value = "  -0.05%  "

# Row 36
$ws.Range("E36").Value = "  -1.15%  "

# Row 37
$ws.Range("E37").Value = "  -0.93%  "

# Row 38
$ws.Range("D38").Value = "'3.33"
$ws.Range("E38").Value = "  -3.27%  "

# Row 39
$ws.Range("E39").Value = "  -0.84%  "

# Row 40
$ws.Range("D40").Value = "'0.994"
$ws.Range("E40").Value = "  -1.04%  "

# Row 41
$ws.Range("D41").Value = "'5.77"
$ws.Range("E41").Value = "  -0.44%  "

# Row 42
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.01%  "

# Row 44
$ws.Range("D44").Value = "'48.10"
$ws.Range("E44").Value = "  +2.06%  "

# Row 45
$ws.Range("D45").Value = "'43.83"
$ws.Range("E45").Value = "  -1.08%  "

# Row 46
$ws.Range("E46").Value = "  -1.56%  "

# Row 47
$ws.Range("D47").Value = "'150.06"
$ws.Range("E47").Value = "  +2.46%  "

# Row 48
$ws.Range("E48").Value = "  -1.71%  "

# Row 49
$ws.Range("D49").Value = "'26.90"
$ws.Range("E49").Value = "  +4.34%  "

# Row 50
$ws.Range("D50").Value = "'389.70"
$ws.Range("E50").Value = "  -1.72%  "

# Row 51
$ws.Range("E51").Value = "  -5.08%  "
